# `tests`: expand excel error sheet
#
# Adds a new row (row 7) to the "cellerrors" sheet:
#   A7 = B7+12            -> #VALUE! (text-formatted, reuses existing style)
#   B7 = C7+20             -> #VALUE! (text-formatted, reuses existing style)
#   C7 = array formula SORT(CHOOSECOLS(A3:B20, 3)) -> #VALUE! (new grey font)
# and moves the sheet's selection to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cellerrors")

# A7 / B7 pick up the existing "text" number-format style (numFmtId 49,
# the same cellXfs entry already used elsewhere in this workbook).
$ws.Range("A7").NumberFormat = "@"
$ws.Range("B7").NumberFormat = "@"

$ws.Range("A7").Formula = "=B7+12"
$ws.Range("B7").Formula = "=C7+20"

# C7 is a dynamic-array formula rendered in a new grey (#454545) font.
$ws.Range("C7").Font.Color = 4539717
$ws.Range("C7").FormulaArray = "=_xlfn._xlws.SORT(_xlfn.CHOOSECOLS(A3:B20, 3))"

# Leave the selection on A8, just past the newly-added row.
$ws.Range("A8").Select() | Out-Null
